$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "60.327.38"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "2.596.85"
$ws.Range("E3").Value = "  -2.89%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.38%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.74%  "
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.37%  "
$ws.Range("D9").Value = "2.604.32"
$ws.Range("E9").Value = "  -2.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.75%  "
$ws.Range("E11").Value = "  -0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.73%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").Value = "3.051.00"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("D15").Value = "60.316.90"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.19%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "2.601.25"
$ws.Range("E18").Value = "  -2.60%  "
$ws.Range("E19").Value = "  -1.43%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "353.51"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.20%  "
$ws.Range("E25").Value = "  -0.59%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.166"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.37%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "0.0₃0839"
$ws.Range("E28").Value = "  -2.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.36"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("E30").Value = "  -0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "19.37"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.72%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "150.92"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.18%  "
$ws.Range("E33").Value = "  -0.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.99"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.55%  "
$ws.Range("E36").Value = "  -2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.878"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.91%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.48"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.841"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.13"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "294.92"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.99%  "
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.620"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("E45").Value = "  -0.23%  "
$ws.Range("E46").Value = "  -4.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.63%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.87"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -1.32%  "
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "1.995.10"
$ws.Range("E51").Value = "  -2.58%  "
